$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 4 (the ConceptScheme row), shifting it down to row 5.
$ws.Rows.Item(4).Insert()

# New row 4: the "tot" Concept entry
$ws.Cells.Item(4, 1).Value = "https://data.bodemenondergrond.vlaanderen.be/id/concept/eenhedenrelatie/tot"
$ws.Cells.Item(4, 2).Value = "http://www.w3.org/2004/02/skos/core#Concept"
$ws.Cells.Item(4, 3).Value = "be.vlaanderen.bodemenondergrond.data.id.concept.eenhedenrelatie.tot"
$ws.Cells.Item(4, 4).Value = "Binnen het aangegeven interval komt een al dan niet volledige opeenvolging van eenheden voor, vanaf de jongste geselecteerde eenheid bovenaan (Lid 1) tot de oudste geselecteerde eenheid onderaan (Lid 2)."
$ws.Cells.Item(4, 5).Value = "https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/eenhedenrelatie"
$ws.Cells.Item(4, 6).Value = "Binnen het aangegeven interval komt een al dan niet volledige opeenvolging van eenheden voor, vanaf de jongste geselecteerde eenheid bovenaan (Lid 1) tot de oudste geselecteerde eenheid onderaan (Lid 2)."
$ws.Cells.Item(4, 7).Value = "tot"
$ws.Cells.Item(4, 8).Value = "https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/eenhedenrelatie"
$ws.Cells.Item(4, 9).Value = "null"

# Row 5 (former row 4, the ConceptScheme entry): update hasTopConcept to include the new "tot" concept
$ws.Cells.Item(5, 9).Value = "https://data.bodemenondergrond.vlaanderen.be/id/concept/eenhedenrelatie/en|https://data.bodemenondergrond.vlaanderen.be/id/concept/eenhedenrelatie/of|https://data.bodemenondergrond.vlaanderen.be/id/concept/eenhedenrelatie/tot"

# Keep the "numbers stored as text" ignored-error marker covering the full used range (A1:I5).
$ws.Range("A1:I5").Errors.Item(513).Ignore = $true
